{"js": "// Replace the 25 division-problem expressions in the table with their\n// updated values (see commit \"Update master to output generated at 503736d\").\n// Each \"from\" string is unique in the document, so a body-wide search for the\n// exact old text and an in-place \"Replace\" insert is safe and unambiguous.\nconst replacements = [\n  [\"90\u00f72=\", \"49\u00f78=\"],\n  [\"69\u00f78=\", \"65\u00f79=\"],\n  [\"16\u00f77=\", \"78\u00f74=\"],\n  [\"49\u00f73=\", \"55\u00f74=\"],\n  [\"84\u00f75=\", \"96\u00f72=\"],\n  [\"77\u00f75=\", \"14\u00f73=\"],\n  [\"26\u00f77=\", \"68\u00f78=\"],\n  [\"24\u00f77=\", \"29\u00f73=\"],\n  [\"58\u00f73=\", \"58\u00f75=\"],\n  [\"88\u00f79=\", \"34\u00f78=\"],\n  [\"30\u00f75=\", \"53\u00f72=\"],\n  [\"99\u00f75=\", \"53\u00f74=\"],\n  [\"14\u00f79=\", \"83\u00f74=\"],\n  [\"39\u00f76=\", \"81\u00f79=\"],\n  [\"67\u00f78=\", \"34\u00f77=\"],\n  [\"72\u00f75=\", \"29\u00f77=\"],\n  [\"71\u00f79=\", \"67\u00f77=\"],\n  [\"80\u00f77=\", \"53\u00f78=\"],\n  [\"33\u00f73=\", \"95\u00f77=\"],\n  [\"32\u00f72=\", \"40\u00f78=\"],\n  [\"47\u00f73=\", \"28\u00f78=\"],\n  [\"63\u00f74=\", \"29\u00f77=\"],\n  [\"84\u00f79=\", \"14\u00f72=\"],\n  [\"62\u00f78=\", \"11\u00f72=\"],\n  [\"51\u00f73=\", \"70\u00f77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Could not find text \"${oldText}\" to replace.`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 division-problem expressions in the table with their\n# updated values (see commit \"Update master to output generated at 503736d\").\n# Each \"from\" string occurs exactly once in the document, so a single\n# Find/Replace pass per pair (wdReplaceOne) is safe and unambiguous.\n$pairs = @(\n  @(\"90\u00f72=\", \"49\u00f78=\"),\n  @(\"69\u00f78=\", \"65\u00f79=\"),\n  @(\"16\u00f77=\", \"78\u00f74=\"),\n  @(\"49\u00f73=\", \"55\u00f74=\"),\n  @(\"84\u00f75=\", \"96\u00f72=\"),\n  @(\"77\u00f75=\", \"14\u00f73=\"),\n  @(\"26\u00f77=\", \"68\u00f78=\"),\n  @(\"24\u00f77=\", \"29\u00f73=\"),\n  @(\"58\u00f73=\", \"58\u00f75=\"),\n  @(\"88\u00f79=\", \"34\u00f78=\"),\n  @(\"30\u00f75=\", \"53\u00f72=\"),\n  @(\"99\u00f75=\", \"53\u00f74=\"),\n  @(\"14\u00f79=\", \"83\u00f74=\"),\n  @(\"39\u00f76=\", \"81\u00f79=\"),\n  @(\"67\u00f78=\", \"34\u00f77=\"),\n  @(\"72\u00f75=\", \"29\u00f77=\"),\n  @(\"71\u00f79=\", \"67\u00f77=\"),\n  @(\"80\u00f77=\", \"53\u00f78=\"),\n  @(\"33\u00f73=\", \"95\u00f77=\"),\n  @(\"32\u00f72=\", \"40\u00f78=\"),\n  @(\"47\u00f73=\", \"28\u00f78=\"),\n  @(\"63\u00f74=\", \"29\u00f77=\"),\n  @(\"84\u00f79=\", \"14\u00f72=\"),\n  @(\"62\u00f78=\", \"11\u00f72=\"),\n  @(\"51\u00f73=\", \"70\u00f77=\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n\n  $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n  if (-not $found) {\n    throw \"Could not find text '$oldText' to replace.\"\n  }\n}\n\nWrite-Output \"done\"\n"}
